$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Header 2 == the "first page" header (word/header1.xml) carrying the
# BTec_Logo-Orange inline picture: rename image1.jpg -> image2.jpg
$hdr = $sec.Headers.Item(2)
if ($hdr.Exists) {
    $ishapes = $hdr.Range.InlineShapes
    for ($i = 1; $i -le $ishapes.Count; $i++) {
        $shp = $ishapes.Item($i)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }
}

# Footer 1 == the primary footer (word/footer1.xml) and
# Footer 2 == the "first page" footer (word/footer2.xml), both carrying
# the Pearson Edexcel logo inline picture: rename image2.png -> image1.png
for ($f = 1; $f -le $sec.Footers.Count; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $ishapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $ishapes.Count; $i++) {
            $shp = $ishapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

Write-Output "done"
